$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do- FY16 Release")

# Insert a new row above row 16, shifting existing rows 16+ down by one.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new to-do item.
$ws.Cells.Item(16, 1).Value = "Not Done"
$ws.Cells.Item(16, 2).Value = "SDK tool Copy-to-clipboard issues due to changes in wxExtGridCtrl"
$ws.Cells.Item(16, 3).Value = "Aron"

# Update the selection to match the post-edit state (A17 selected).
$ws.Range("A17").Select()
